$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS group): VENTA and POR CUMPLIR updated
$ws.Range("D2").Value = 1184.9
$ws.Range("E2").Value = -1184.9

# Row 4 (TOTAL row): VENTA, POR CUMPLIR and CUMPLIMIENTO updated
$ws.Range("D4").Value = 1737.33
$ws.Range("E4").Value = 11986.01
$ws.Range("F4").Value = 0.1265967322823744
